$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2675185060000028
$ws.Range("H2").Value = 0.987

$ws.Range("G3").Value = 0.2675185060000028
$ws.Range("H3").Value = 0.987
